# Juno: check in to OLPRODLOC.
# Localizes the "Sales report" worksheet to Arabic: the sheet itself is
# renamed and the column headers (row 1) are translated. The quarterly
# labels in column A (2022-Q1 .. 2023-Q4) and all numeric sales figures
# are left exactly as they were.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet: "Sales report" -> "تقرير المبيعات"
$ws.Name = "تقرير المبيعات"

# Translate the header row (A1:G1) into Arabic.
$ws.Range("A1").Value = "ربع السنوي"        # Year-Quarter
$ws.Range("B1").Value = "الغرب الأوسط"      # Midwest
$ws.Range("C1").Value = "جبل"               # Mountain
$ws.Range("D1").Value = "شمال شرق"          # Northeast
$ws.Range("E1").Value = "الجنوب"            # South
$ws.Range("F1").Value = "جنوب شرق"          # Southeast
$ws.Range("G1").Value = "الغرب"             # West

# Re-select the worksheet so it stays the active tab of the workbook.
$ws.Activate()
